{"js": "// Applies the \"feat: work on pres + report\" edits to the Elixir report:\n//  - justify (w:jc=\"both\") a handful of paragraphs\n//  - duplicate one empty paragraph (Documentation section)\n//  - a batch of small wording fixes (grammar / typos / clarifications)\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Helper: find the unique range matching `query` and replace its text.\n// ---------------------------------------------------------------------\nasync function replaceText(query, replacement) {\n  const results = body.search(query, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"replaceText: no match for \" + JSON.stringify(query));\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Helper: find the unique range matching `query` and return its first\n// paragraph (after syncing the search).\nasync function paragraphFor(query) {\n  const results = body.search(query, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"paragraphFor: no match for \" + JSON.stringify(query));\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  await context.sync();\n  return para;\n}\n\n// ---------------------------------------------------------------------\n// 1) Justify the intro paragraph (\"Elixir a \u00e9t\u00e9 cr\u00e9\u00e9 par Jos\u00e9 Valim...\")\n// ---------------------------------------------------------------------\nconst introPara = await paragraphFor(\"Elixir a \u00e9t\u00e9 cr\u00e9\u00e9 par Jos\u00e9 Valim\");\nintroPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) \"Cette phrase, bien non sans int\u00e9r\u00eat,\" -> \"... bien que non sans int\u00e9r\u00eat,\"\n// ---------------------------------------------------------------------\nawait replaceText(\" bien non sans int\u00e9r\u00eat, \", \" bien que non sans int\u00e9r\u00eat, \");\n\n// ---------------------------------------------------------------------\n// 3) \"ex\u00e9cuter\" -> \"ex\u00e9cut\u00e9\" (wrong infinitive -> correct past participle)\n// ---------------------------------------------------------------------\nawait replaceText(\n  \"L\u2019int\u00e9gralit\u00e9 du code Elixir est ex\u00e9cuter dans des \",\n  \"L\u2019int\u00e9gralit\u00e9 du code Elixir est ex\u00e9cut\u00e9 dans des \"\n);\n\n// ---------------------------------------------------------------------\n// 4) \"confondre processus Elixir\" -> \"confondre les processus Elixir\"\n// ---------------------------------------------------------------------\nawait replaceText(\n  \". Il ne faut pas confondre processus Elixir\",\n  \". Il ne faut pas confondre les processus Elixir\"\n);\n\n// ---------------------------------------------------------------------\n// 5) Justify the \"Documentation\" paragraph + surrounding empty paragraphs,\n//    duplicate the empty paragraph before \"De plus, ...\", and fix\n//    \"toute \u00e9l\u00e9ment\" -> \"tout \u00e9l\u00e9ment\".\n// ---------------------------------------------------------------------\nconst docPara = await paragraphFor(\"Elixir traite la documentation comme\");\ndocPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\nconst deplusPara = await paragraphFor(\"De plus, toute \u00e9l\u00e9ment du langage de base\");\nconst emptyPara = deplusPara.getPrevious();\nemptyPara.load(\"text\");\nawait context.sync();\n\nemptyPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\n// Duplicate that empty paragraph (the diff shows a brand-new empty\n// paragraph inserted right before the \"De plus\" one).\nconst duplicatedEmptyPara = emptyPara.insertParagraph(\"\", Word.InsertLocation.after);\nduplicatedEmptyPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\ndeplusPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\nawait replaceText(\n  \"De plus, toute \u00e9l\u00e9ment du langage de base\",\n  \"De plus, tout \u00e9l\u00e9ment du langage de base\"\n);\n\n// ---------------------------------------------------------------------\n// 6) Justify the \"Les atomes...\" paragraph, add \"le caract\u00e8re \u00ab : \u00bb\"\n//    clarification and split the \"nombres ... compr\u00e9hensibles\" sentence.\n// ---------------------------------------------------------------------\nconst atomsPara = await paragraphFor(\"Les atomes sont des constantes avec comme valeur\");\natomsPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\nawait replaceText(\n  \"doivent commencer par : et sont grandement utilis\u00e9s\",\n  \"doivent commencer par le caract\u00e8re \u00ab : \u00bb et sont grandement utilis\u00e9s\"\n);\n\nawait replaceText(\n  \"cha\u00eenes de caract\u00e8res en dur ou de nombres peu compr\u00e9hensibles.\",\n  \"cha\u00eenes de caract\u00e8res en dur ou de nombres indiquant des \u00e9tats peu compr\u00e9hensibles.\"\n);\n\n// ---------------------------------------------------------------------\n// 7) \"un op\u00e9rateur simple d\u2019affectation\" -> \"un simple op\u00e9rateur d\u2019affectation\"\n// ---------------------------------------------------------------------\nawait replaceText(\n  \"Pour l\u2019anecdote, l\u2019op\u00e9rateur = n\u2019est pas un op\u00e9rateur simple d\u2019affectation en Elixir\",\n  \"Pour l\u2019anecdote, l\u2019op\u00e9rateur = n\u2019est pas un simple op\u00e9rateur d\u2019affectation en Elixir\"\n);\n\n// ---------------------------------------------------------------------\n// 8) \"temps r\u00e9els\" -> \"temps r\u00e9el\"\n// ---------------------------------------------------------------------\nawait replaceText(\n  \"synchroniser en temps r\u00e9els tous les utilisateurs.\",\n  \"synchroniser en temps r\u00e9el tous les utilisateurs.\"\n);\n", "ps1": "# Applies the \"feat: work on pres + report\" edits to the Elixir report:\n#  - justify (w:jc=\"both\") a handful of paragraphs\n#  - duplicate one empty paragraph (Documentation section)\n#  - a batch of small wording fixes (grammar / typos / clarifications)\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2, wdFindContinue = 1, wdAlignParagraphJustify = 3\n$wdReplaceAll = 2\n$wdFindContinue = 1\n$wdAlignParagraphJustify = 3\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\nfunction Get-ParagraphContaining($searchText) {\n    $rng = $d.Content\n    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false)\n    return $rng.Paragraphs(1)\n}\n\n# ---------------------------------------------------------------------\n# 1) Justify the intro paragraph (\"Elixir a \u00e9t\u00e9 cr\u00e9\u00e9 par Jos\u00e9 Valim...\")\n# ---------------------------------------------------------------------\n$introPara = Get-ParagraphContaining(\"Elixir a \u00e9t\u00e9 cr\u00e9\u00e9 par Jos\u00e9 Valim\")\n$introPara.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify\n\n# ---------------------------------------------------------------------\n# 2) \"Cette phrase, bien non sans int\u00e9r\u00eat,\" -> \"... bien que non sans int\u00e9r\u00eat,\"\n# ---------------------------------------------------------------------\nReplace-Text \" bien non sans int\u00e9r\u00eat, \" \" bien que non sans int\u00e9r\u00eat, \"\n\n# ---------------------------------------------------------------------\n# 3) \"ex\u00e9cuter\" -> \"ex\u00e9cut\u00e9\" (wrong infinitive -> correct past participle)\n# ---------------------------------------------------------------------\nReplace-Text \"L\u2019int\u00e9gralit\u00e9 du code Elixir est ex\u00e9cuter dans des \" \"L\u2019int\u00e9gralit\u00e9 du code Elixir est ex\u00e9cut\u00e9 dans des \"\n\n# ---------------------------------------------------------------------\n# 4) \"confondre processus Elixir\" -> \"confondre les processus Elixir\"\n# ---------------------------------------------------------------------\nReplace-Text \". Il ne faut pas confondre processus Elixir\" \". Il ne faut pas confondre les processus Elixir\"\n\n# ---------------------------------------------------------------------\n# 5) Justify the \"Documentation\" paragraph + surrounding empty paragraphs,\n#    duplicate the empty paragraph before \"De plus, ...\", and fix\n#    \"toute \u00e9l\u00e9ment\" -> \"tout \u00e9l\u00e9ment\".\n# ---------------------------------------------------------------------\n$docPara = Get-ParagraphContaining(\"Elixir traite la documentation comme\")\n$docPara.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify\n\n$deplusPara = Get-ParagraphContaining(\"De plus, toute \u00e9l\u00e9ment du langage de base\")\n$emptyPara = $deplusPara.Previous()\n$emptyPara.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify\n\n# Duplicate that empty paragraph (the diff shows a brand-new empty\n# paragraph inserted right before the \"De plus\" one).\n$emptyPara.Range.InsertParagraphAfter()\n\n# Re-fetch the \"De plus\" paragraph and justify it along with the newly\n# inserted empty paragraph (insertion shifted paragraph indices).\n$deplusPara = Get-ParagraphContaining(\"De plus, toute \u00e9l\u00e9ment du langage de base\")\n$newEmptyPara = $deplusPara.Previous()\n$newEmptyPara.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify\n$deplusPara.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify\n\nReplace-Text \"De plus, toute \u00e9l\u00e9ment du langage de base\" \"De plus, tout \u00e9l\u00e9ment du langage de base\"\n\n# ---------------------------------------------------------------------\n# 6) Justify the \"Les atomes...\" paragraph, add \"le caract\u00e8re \u00ab : \u00bb\"\n#    clarification and split the \"nombres ... compr\u00e9hensibles\" sentence.\n# ---------------------------------------------------------------------\n$atomsPara = Get-ParagraphContaining(\"Les atomes sont des constantes avec comme valeur\")\n$atomsPara.Range.ParagraphFormat.Alignment = $wdAlignParagraphJustify\n\nReplace-Text \"doivent commencer par : et sont grandement utilis\u00e9s\" \"doivent commencer par le caract\u00e8re \u00ab : \u00bb et sont grandement utilis\u00e9s\"\n\nReplace-Text \"cha\u00eenes de caract\u00e8res en dur ou de nombres peu compr\u00e9hensibles.\" \"cha\u00eenes de caract\u00e8res en dur ou de nombres indiquant des \u00e9tats peu compr\u00e9hensibles.\"\n\n# ---------------------------------------------------------------------\n# 7) \"un op\u00e9rateur simple d'affectation\" -> \"un simple op\u00e9rateur d'affectation\"\n# ---------------------------------------------------------------------\nReplace-Text \"Pour l\u2019anecdote, l\u2019op\u00e9rateur = n\u2019est pas un op\u00e9rateur simple d\u2019affectation en Elixir\" \"Pour l\u2019anecdote, l\u2019op\u00e9rateur = n\u2019est pas un simple op\u00e9rateur d\u2019affectation en Elixir\"\n\n# ---------------------------------------------------------------------\n# 8) \"temps r\u00e9els\" -> \"temps r\u00e9el\"\n# ---------------------------------------------------------------------\nReplace-Text \"synchroniser en temps r\u00e9els tous les utilisateurs.\" \"synchroniser en temps r\u00e9el tous les utilisateurs.\"\n"}
